# Apply the target changes to WS_holdings workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclosure text (cell A16) with new as-of date
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-21 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-13
$ws.Range("D2").Value = 0.03096465881679771
$ws.Range("E2").Value = 0.006887052341597588

$ws.Range("D3").Value = 0.02463517182500707
$ws.Range("E3").Value = -0.006400000000000072

$ws.Range("D4").Value = 0.05203379592145786
$ws.Range("E4").Value = 0.005401596993893865

$ws.Range("D5").Value = 0.1386681324568943
$ws.Range("E5").Value = 0.01420640104506843

$ws.Range("D6").Value = 0.02840799223080062
$ws.Range("E6").Value = 0.01498422712933767

$ws.Range("D7").Value = 0.1206733449896883
$ws.Range("E7").Value = 0.01374353801538275

$ws.Range("D8").Value = 0.09992828394478079
$ws.Range("E8").Value = 0.01439666603523393

$ws.Range("D9").Value = 0.02777661299423995
$ws.Range("E9").Value = 0.01869775626924763

$ws.Range("D10").Value = 0.1190874689072682
$ws.Range("E10").Value = 0.01488933601609665

$ws.Range("D11").Value = 0.2549933471557863
$ws.Range("E11").Value = 0.01026602392429932

$ws.Range("D12").Value = 0.1028311907572789
$ws.Range("E12").Value = 0.003131728322568028

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.01106171630202013

# Restore sheet protection (original state before the edit)
$ws.Protect()
